$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2372
$ws.Range("I43").Value = 1215
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 1215
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = -1146
$ws.Range("N43").Value = -7138
$ws.Range("H113").Value = 10153.846
$ws.Range("J113").Value = 19500
$ws.Range("L113").Value = 19500
$ws.Range("N113").Value = -26008
$ws.Range("H116").Value = 349666.78
$ws.Range("I116").Value = 715674.9399999999
$ws.Range("J116").Value = 8059.2
$ws.Range("K116").Value = 715674.9399999999
$ws.Range("L116").Value = 8059.2
$ws.Range("M116").Value = -712232.9399999999
$ws.Range("N116").Value = -14943.2
$ws.Range("H118").Value = 1525
$ws.Range("I118").Value = 1634
$ws.Range("K118").Value = 4902
$ws.Range("M118").Value = -3245
$ws.Range("H132").Value = 29826852
$ws.Range("I132").Value = 39001810
$ws.Range("J132").Value = 8248.5
$ws.Range("K132").Value = 117005430
$ws.Range("L132").Value = 24745.5
$ws.Range("M132").Value = -117002900
$ws.Range("N132").Value = -29805.5
$ws.Range("H139").Value = 48884
$ws.Range("J139").Value = 48884
$ws.Range("L139").Value = 48884
$ws.Range("N139").Value = -59164

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4903.6323
$ws.Range("I32").Value = 5029.0815
$ws.Range("J32").Value = 4580.1055
$ws.Range("K32").Value = 5029.0815
$ws.Range("L32").Value = 4580.1055
$ws.Range("M32").Value = -4742.0815
$ws.Range("N32").Value = -5154.1055
$ws.Range("H45").Value = 2607.111
$ws.Range("I45").Value = 1940.6923
$ws.Range("J45").Value = 4339.8
$ws.Range("K45").Value = 1940.6923
$ws.Range("L45").Value = 4339.8
$ws.Range("M45").Value = -1563.6923
$ws.Range("N45").Value = -5093.8
$ws.Range("H109").Value = 35500
$ws.Range("J109").Value = 35500
$ws.Range("L109").Value = 35500
$ws.Range("N109").Value = -38274
$ws.Range("H122").Value = 2046.5714
$ws.Range("I122").Value = 1378.3
$ws.Range("J122").Value = 3717.25
$ws.Range("K122").Value = 4134.9
$ws.Range("L122").Value = 11151.75
$ws.Range("M122").Value = -1684.9
$ws.Range("N122").Value = -16051.75
$ws.Range("H132").Value = 2844.5144
$ws.Range("I132").Value = 2011.75
$ws.Range("J132").Value = 3954.8667
$ws.Range("K132").Value = 6035.25
$ws.Range("L132").Value = 11864.6001
$ws.Range("M132").Value = -3505.25
$ws.Range("N132").Value = -16924.6001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 46572.633
$ws.Range("J140").Value = 46572.633
$ws.Range("L140").Value = 46572.633
$ws.Range("N140").Value = -56932.633

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2209.027
$ws.Range("I122").Value = 1690.16
$ws.Range("J122").Value = 3290
$ws.Range("K122").Value = 5070.48
$ws.Range("L122").Value = 9870
$ws.Range("M122").Value = -2620.48
$ws.Range("N122").Value = -14770
$ws.Range("H132").Value = 3105.147
$ws.Range("I132").Value = 2845.7
$ws.Range("J132").Value = 3475.7856
$ws.Range("K132").Value = 8537.099999999999
$ws.Range("L132").Value = 10427.3568
$ws.Range("M132").Value = -6007.099999999999
$ws.Range("N132").Value = -15487.3568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2910.984
$ws.Range("J68").Value = 3166.0378
$ws.Range("L68").Value = 9498.1134
$ws.Range("N68").Value = -11120.1134
$ws.Range("H71").Value = 2910.984
$ws.Range("J71").Value = 3166.0378
$ws.Range("L71").Value = 28494.3402
$ws.Range("N71").Value = -36606.3402
$ws.Range("H113").Value = 457.59616
$ws.Range("I113").Value = 454.7353
$ws.Range("K113").Value = 1364.2059
$ws.Range("M113").Value = 805.7941000000001
$ws.Range("H123").Value = 9512.5
$ws.Range("I123").Value = 4900
$ws.Range("J123").Value = 11050
$ws.Range("K123").Value = 14700
$ws.Range("L123").Value = 33150
$ws.Range("M123").Value = -12250
$ws.Range("N123").Value = -38050
$ws.Range("H125").Value = 7000
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 30000
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -39840
$ws.Range("H131").Value = 779.58
$ws.Range("J131").Value = 806.15216
$ws.Range("L131").Value = 2418.45648
$ws.Range("N131").Value = -12498.45648
$ws.Range("H137").Value = 3441.611
$ws.Range("I137").Value = 3212.0715
$ws.Range("K137").Value = 9636.2145
$ws.Range("M137").Value = -4536.2145
$ws.Range("H140").Value = 983.8095
$ws.Range("I140").Value = 547.7778
$ws.Range("J140").Value = 3600
$ws.Range("K140").Value = 1643.3334
$ws.Range("L140").Value = 10800
$ws.Range("M140").Value = 3536.6666
$ws.Range("N140").Value = -21160

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6173701.5
$ws.Range("I107").Value = 523.55554
$ws.Range("J107").Value = 12346879
$ws.Range("K107").Value = 523.55554
$ws.Range("L107").Value = 12346879
$ws.Range("M107").Value = 1396.44446
$ws.Range("N107").Value = -12350719
$ws.Range("H122").Value = 2834.9092
$ws.Range("I122").Value = 1705.4
$ws.Range("J122").Value = 3776.1667
$ws.Range("K122").Value = 5116.200000000001
$ws.Range("L122").Value = 11328.5001
$ws.Range("M122").Value = -2666.200000000001
$ws.Range("N122").Value = -16228.5001
$ws.Range("H123").Value = 10509.857
$ws.Range("J123").Value = 10509.857
$ws.Range("L123").Value = 10509.857
$ws.Range("N123").Value = -15409.857
$ws.Range("H126").Value = 3786.7761
$ws.Range("I126").Value = 2946.8086
$ws.Range("J126").Value = 5760.7
$ws.Range("K126").Value = 8840.425799999999
$ws.Range("L126").Value = 17282.1
$ws.Range("M126").Value = -6370.425799999999
$ws.Range("N126").Value = -22222.1
$ws.Range("H132").Value = 4335.973
$ws.Range("I132").Value = 3623.611
$ws.Range("J132").Value = 5010.8423
$ws.Range("K132").Value = 10870.833
$ws.Range("L132").Value = 15032.5269
$ws.Range("M132").Value = -8340.832999999999
$ws.Range("N132").Value = -20092.5269

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3776.111
$ws.Range("I7").Value = 1796
$ws.Range("J7").Value = 6251.25
$ws.Range("K7").Value = 1796
$ws.Range("L7").Value = 6251.25
$ws.Range("M7").Value = -1684
$ws.Range("N7").Value = -6475.25
$ws.Range("H122").Value = 3849.8
$ws.Range("I122").Value = 2970.5264
$ws.Range("J122").Value = 6634.1665
$ws.Range("K122").Value = 8911.5792
$ws.Range("L122").Value = 19902.4995
$ws.Range("M122").Value = -6461.5792
$ws.Range("N122").Value = -24802.4995
$ws.Range("H126").Value = 3776.111
$ws.Range("I126").Value = 1796
$ws.Range("J126").Value = 6251.25
$ws.Range("K126").Value = 5388
$ws.Range("L126").Value = 18753.75
$ws.Range("M126").Value = -2918
$ws.Range("N126").Value = -23693.75
$ws.Range("H132").Value = 3765.9614
$ws.Range("I132").Value = 2784.2354
$ws.Range("J132").Value = 5620.3335
$ws.Range("K132").Value = 8352.706200000001
$ws.Range("L132").Value = 16861.0005
$ws.Range("M132").Value = -5822.706200000001
$ws.Range("N132").Value = -21921.0005
$ws.Range("H133").Value = 30028.75
$ws.Range("J133").Value = 30028.75
$ws.Range("L133").Value = 30028.75
$ws.Range("N133").Value = -35088.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 73130060
$ws.Range("I96").Value = 90909810
$ws.Range("J96").Value = 7937661.5
$ws.Range("K96").Value = 90909810
$ws.Range("L96").Value = 7937661.5
$ws.Range("M96").Value = -90908437
$ws.Range("N96").Value = -7940407.5
$ws.Range("H122").Value = 2848.5
$ws.Range("I122").Value = 980.9091
$ws.Range("J122").Value = 5131.1113
$ws.Range("K122").Value = 2942.7273
$ws.Range("L122").Value = 15393.3339
$ws.Range("M122").Value = -492.7273
$ws.Range("N122").Value = -20293.3339
